# Survival.metrics IQ document — apply the "Add Validation procedure to
# Export.metrics" edit set.
#
# wdReplace constants used with Find.Execute(..., Replace:=N):
#   0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll
# wdFindWrap:
#   1 = wdFindContinue

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found -> $old"
    }
}

# 1) Merge the three runs that spell out the package name in the purpose
#    paragraph into a single run (no textual change, just a run-merge that
#    happens naturally once the range is rewritten).
Replace-Text "본 문서의 목적은 R 패키지 Survival.metrics(v0.1.0) 및 관련 분석 환경이 사전 정의된 사양(" `
             "본 문서의 목적은 R 패키지 Survival.metrics(v0.1.0) 및 관련 분석 환경이 사전 정의된 사양("

# 2) R version bump 4.5.1 -> 4.5.2
Replace-Text "R version: 4.5.1" "R version: 4.5.2"

# 3) Evidence-collection bullet: log file now carries the package prefix
Replace-Text " 3단계와 4단계를 실행하는 과정에서 R 콘솔에 출력되는 모든 텍스트 로그를 복사하여, 검증 증거 자료로서 install_log.txt 파일로 저장" `
             " 3단계와 4단계를 실행하는 과정에서 R 콘솔에 출력되는 모든 텍스트 로그를 복사하여, 검증 증거 자료로서 Survival.metrics_install_log.txt 파일로 저장"

# 4) Acceptance-criteria bullet: same rename
Replace-Text " 5단계에 따라 생성된 install_log.txt 파일이 존재해야 하며, 여기에는 " `
             " 5단계에 따라 생성된 Survival.metrics_install_log.txt 파일이 존재해야 하며, 여기에는 "

# 5) Merge the date runs (2025년 10월 31일) into a single run (no text change)
Replace-Text "본 IQ는 2025년 10월 31일에 " "본 IQ는 2025년 10월 31일에 "

# 6) Merge the runs referencing the install log filename at the end of the
#    "실행 결과 및 증거" section (no text change)
Replace-Text " 실행 시 콘솔에 출력된 전체 로그는 본 보고서와 함께 제출된 Survival.metrics_install_log.txt 파일을 참조함." `
             " 실행 시 콘솔에 출력된 전체 로그는 본 보고서와 함께 제출된 Survival.metrics_install_log.txt 파일을 참조함."

# 8) Validation-criteria-summary table, first data row / first column:
#    rename the referenced log file there too.
Replace-Text "install_log.txt 파일이 존재해야 함" "Survival.metrics_install_log.txt 파일이 존재해야 함"

# 7) The table's nominal grid widths shift slightly now that the first
#    column holds more text (column totals are unchanged: 3574+4168 ==
#    3611+4131 == 7742 twips); set the two affected columns explicitly.
$summaryTable = $d.Tables(1)
$summaryTable.Columns(1).Width = 180.55   # 3611 twips
$summaryTable.Columns(2).Width = 206.55   # 4131 twips

# 9) The trailing empty paragraph (right before the final sectPr) loses its
#    stray direct formatting (an East-Asian font hint with no visible
#    effect), becoming a fully bare paragraph mark.
$last = $d.Paragraphs.Last
$last.Range.Select()
$word.Selection.ClearFormatting()

Write-Output "done"
